# The workbook contains a weekly price log for "Hortaliza, Agrícola del Norte
# S.A. de Arica - Pimiento". A new week of data (date serial 44641, i.e.
# 2022-03-21) is inserted at the top of the data block (rows 433-438),
# pushing all the existing historical rows down by 6 rows (old row 433
# becomes row 439, ... old row 543 becomes row 549).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows starting at row 433 - this shifts every row from 433
# downward by 6 rows and extends the sheet's used range / dimension
# automatically (A1:R543 -> A1:R549).
$ws.Range("A433:A438").EntireRow.Insert()

# Values that are constant across every data row in this sheet.
$mercadoId = 1
$mercado = "Agrícola del Norte S.A. de Arica"
$region = "Arica y Parinacota"
$codreg = 15
$categoriaId = 100112002
$categoria = "Pimiento"
$unidad = "`$/caja 15 kilos"
$origen = "Región de Arica y Parinacota"
$kgUnidades = 15
$clasificacion = "Hortaliza"

# New data for the inserted rows: date, variety, quality, volume,
# min price, max price, weighted avg price, price/kg.
$fecha = 44641

$newRows = @(
    @{ Row = 433; Variedad = "Zafiro rojo";  Calidad = "Primera"; Volumen = 100; PMin = 24000; PMax = 25000; PProm = 24500; PKg = 1633 },
    @{ Row = 434; Variedad = "Zafiro rojo";  Calidad = "Segunda"; Volumen = 140; PMin = 21000; PMax = 22000; PProm = 21500; PKg = 1433 },
    @{ Row = 435; Variedad = "Zafiro rojo";  Calidad = "Tercera"; Volumen = 150; PMin = 18000; PMax = 20000; PProm = 19000; PKg = 1267 },
    @{ Row = 436; Variedad = "Zafiro verde"; Calidad = "Primera"; Volumen = 120; PMin = 13000; PMax = 14000; PProm = 13500; PKg = 900 },
    @{ Row = 437; Variedad = "Zafiro verde"; Calidad = "Segunda"; Volumen = 130; PMin = 11000; PMax = 12000; PProm = 11500; PKg = 767 },
    @{ Row = 438; Variedad = "Zafiro verde"; Calidad = "Tercera"; Volumen = 140; PMin = 9000;  PMax = 10000; PProm = 9500;  PKg = 633 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $categoriaId
    $ws.Cells.Item($r, 7).Value = $categoria
    $ws.Cells.Item($r, 8).Value = $item.Variedad
    $ws.Cells.Item($r, 9).Value = $item.Calidad
    $ws.Cells.Item($r, 10).Value = $item.Volumen
    $ws.Cells.Item($r, 11).Value = $item.PMin
    $ws.Cells.Item($r, 12).Value = $item.PMax
    $ws.Cells.Item($r, 13).Value = $item.PProm
    $ws.Cells.Item($r, 14).Value = $unidad
    $ws.Cells.Item($r, 15).Value = $origen
    $ws.Cells.Item($r, 16).Value = $item.PKg
    $ws.Cells.Item($r, 17).Value = $kgUnidades
    $ws.Cells.Item($r, 18).Value = $clasificacion
}
